# Refresh the cryptos list prices / 1h volume figures (and restore the
# Polkadot / WrappedEther row order) as captured in the source diff.
# Numeric-looking text (prices, percentages) is written with a leading
# apostrophe so Excel keeps it as literal text (e.g. '591.00', '2.706.24',
# "  +0.18%  ") instead of silently coercing it to a number and losing
# formatting such as trailing zeros / thousands separators / padding.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.000.51'
$ws.Range('E2').Value = '''  -0.01%  '
$ws.Range('D3').Value = '''2.921.68'
$ws.Range('E3').Value = '''  -0.08%  '
$ws.Range('E4').Value = '''  +0.01%  '
$ws.Range('D5').Value = '''591.00'
$ws.Range('E5').Value = '''  +0.61%  '
$ws.Range('D6').Value = '''146.86'
$ws.Range('E6').Value = '''  +0.48%  '
$ws.Range('E7').Value = '''  +0.04%  '
$ws.Range('E8').Value = '''  +0.18%  '
$ws.Range('D9').Value = '''6.92'
$ws.Range('E9').Value = '''  -0.36%  '
$ws.Range('D10').Value = '''0.143'
$ws.Range('E10').Value = '''  -1.08%  '
$ws.Range('E11').Value = '''  -1.76%  '
$ws.Range('E12').Value = '''  -0.19%  '
$ws.Range('D13').Value = '''33.65'
$ws.Range('E13').Value = '''  -0.17%  '
$ws.Range('E14').Value = '''  -0.16%  '
$ws.Range('D15').Value = '''3.406.56'
$ws.Range('E15').Value = '''  -0.01%  '
$ws.Range('D16').Value = '''60.958.97'
$ws.Range('E16').Value = '''  +0.01%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''6.70'
$ws.Range('E17').Value = '''  -1.21%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '''2.922.18'
$ws.Range('E18').Value = '''  +0.00%  '
$ws.Range('D19').Value = '''432.82'
$ws.Range('E19').Value = '''  +0.42%  '
$ws.Range('E20').Value = '''  -1.62%  '
$ws.Range('E21').Value = '''  -0.77%  '
$ws.Range('D22').Value = '''7.12'
$ws.Range('E22').Value = '''  -0.39%  '
$ws.Range('D23').Value = '''81.32'
$ws.Range('E23').Value = '''  +0.91%  '
$ws.Range('D24').Value = '''10.91'
$ws.Range('E24').Value = '''  +0.28%  '
$ws.Range('D25').Value = '''2.22'
$ws.Range('E25').Value = '''  -0.88%  '
$ws.Range('D26').Value = '''11.88'
$ws.Range('E26').Value = '''  -0.66%  '
$ws.Range('E27').Value = '''  -0.04%  '
$ws.Range('D28').Value = '''2.26'
$ws.Range('E28').Value = '''  +4.19%  '
$ws.Range('E29').Value = '''  -0.48%  '
$ws.Range('E30').Value = '''  -3.62%  '
$ws.Range('E31').Value = '''  +0.23%  '
$ws.Range('E32').Value = '''  +2.66%  '
$ws.Range('E33').Value = '''  +0.03%  '
$ws.Range('D34').Value = '0.0₃0868'
$ws.Range('E34').Value = '''  -0.99%  '
$ws.Range('E35').Value = '''  -0.21%  '
$ws.Range('E36').Value = '''  -0.34%  '
$ws.Range('E37').Value = '''  -0.89%  '
$ws.Range('E38').Value = '''  -1.81%  '
$ws.Range('D39').Value = '''0.122'
$ws.Range('E39').Value = '''  -5.43%  '
$ws.Range('E40').Value = '''  -1.53%  '
$ws.Range('D41').Value = '''41.67'
$ws.Range('E41').Value = '''  +0.16%  '
$ws.Range('E42').Value = '''  -5.41%  '
$ws.Range('D43').Value = '''377.61'
$ws.Range('E43').Value = '''  -0.54%  '
$ws.Range('E44').Value = '''  -1.41%  '
$ws.Range('D45').Value = '''2.706.24'
$ws.Range('E45').Value = '''  +0.18%  '
$ws.Range('D46').Value = '''133.77'
$ws.Range('E46').Value = '''  +0.95%  '
$ws.Range('D48').Value = '''23.92'
$ws.Range('E48').Value = '''  -4.43%  '
$ws.Range('E49').Value = '''  -0.70%  '
$ws.Range('E50').Value = '''  -3.02%  '
$ws.Range('E51').Value = '''  -0.89%  '
